$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 62 ("Gene view" section header),
# shifting rows 62-90 down to 63-91.
$ws.Rows.Item(62).Insert() | Out-Null

# Populate the newly inserted row with the new test-step text describing
# clicking a patient's name next to the bar, and the resulting behaviour.
$ws.Range("A62").Value = "Click on a name of a patient next to the bar"
$ws.Range("B62").Value = "A div should appear with the whole phenotype of the patient. This div is resizable and draggable."

# The inserted row copied formatting (incl. wrap text) from the row above;
# match the plain (non-wrapped) style used by the rest of column A/B cells.
$ws.Range("A62").Copy() | Out-Null
$ws.Range("B62").PasteSpecial(-4122) | Out-Null

# Update the active selection to match the edited area.
$ws.Range("A13").Select() | Out-Null
